$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old lingering cell selection (C12) left over from editing the
# table, returning the cursor to the top-left cell.
$ws.Range("A1").Select() | Out-Null

# Row 20 previously held the "Azfar & Nelson" lab estimate (with several blank
# cells). That row is replaced with the (previously row 21) "Solaz, De Vries,
# & de Geus" lab estimate, now fully populated.
$ws.Range("B20").Value = 2018
$ws.Range("C20").Value = "Solaz, De Vries, & de Geus"
$ws.Range("D20").Value = "Solaz et al."
$ws.Range("E20").Value = "UK"
$ws.Range("F20").Value = -0.392
$ws.Range("G20").Value = 0.163
$ws.Range("J20").Value = 1168

# Row 21 becomes a brand-new lab experiment point estimate: Arvate & Mittlaender.
$ws.Range("B21").Value = 2017
$ws.Range("C21").Value = "Arvate & Mittlaender"
$ws.Range("D21").Value = "Arvate & Mittlaender"
$ws.Range("E21").Value = "Brazil"
$ws.Range("F21").Value = -0.46
$ws.Range("G21").Value = 0.16
$ws.Range("J21").Value = 80
$ws.Range("L21").Value = "N is number of respondents. Ten rounds. SE clustered by respondent. "
